$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").Value = 46027
$ws.Range("A78").NumberFormat = $ws.Range("A77").NumberFormat

$ws.Range("B78").Value = 175
$ws.Range("C78").Value = 182
$ws.Range("D78").Value = 178
